$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing changed values (recalculated stats) ---
$ws.Range("G3").Value = 1.77119469641527
$ws.Range("G4").Value = 0.0119399594285355
$ws.Range("G5").Value = 0.0119399594285355
$ws.Range("F11").Value = 0.00208
$ws.Range("G11").Value = 0.0079508567877498
$ws.Range("L11").Value = 0.00208
$ws.Range("M11").Value = 0.01414
$ws.Range("F12").Value = 0.00208
$ws.Range("G12").Value = 0.0079508567877498
$ws.Range("L12").Value = 0.00208
$ws.Range("M12").Value = 0.01414
$ws.Range("G13").Value = 0.491883212265945
$ws.Range("G14").Value = 0.491883212265945
$ws.Range("F16").Value = 0.3782
$ws.Range("G16").Value = 0.50683559234422
$ws.Range("L16").Value = 0.08805
$ws.Range("M16").Value = 1.02166
$ws.Range("F17").Value = 0.3782
$ws.Range("G17").Value = 0.50683559234422
$ws.Range("L17").Value = 0.08805
$ws.Range("M17").Value = 1.02166
$ws.Range("G23").Value = 1.6132522180999
$ws.Range("G24").Value = 0.0126189431670349
$ws.Range("G25").Value = 0.0126189431670349
$ws.Range("G26").Value = 1407.74990941775
$ws.Range("H26").Value = 13516.3435359425
$ws.Range("G27").Value = 1407.74990941775
$ws.Range("H27").Value = 13516.3435359425
$ws.Range("G28").Value = 1407.74990941775
$ws.Range("H28").Value = 13516.3435359425
$ws.Range("G29").Value = 1407.74990941775
$ws.Range("H29").Value = 13516.3435359425
$ws.Range("F31").Value = 0.00276
$ws.Range("G31").Value = 0.0100512651519251
$ws.Range("M31").Value = 0.01663
$ws.Range("F32").Value = 0.00276
$ws.Range("G32").Value = 0.0100512651519251
$ws.Range("M32").Value = 0.01663
$ws.Range("G33").Value = 0.516565739517869
$ws.Range("G34").Value = 0.516565739517869
$ws.Range("G36").Value = 0.533006076619067
$ws.Range("G37").Value = 0.533006076619067
$ws.Range("G44").Value = 0.0118365111526735
$ws.Range("G45").Value = 0.0118365111526735
$ws.Range("G46").Value = 1061.39025711662
$ws.Range("H46").Value = 13516.3435359425
$ws.Range("G47").Value = 1061.39025711662
$ws.Range("H47").Value = 13516.3435359425
$ws.Range("G48").Value = 1061.39025711662
$ws.Range("H48").Value = 13516.3435359425
$ws.Range("G49").Value = 1061.39025711662
$ws.Range("H49").Value = 13516.3435359425
$ws.Range("F51").Value = 0.00404
$ws.Range("G51").Value = 0.0118370784417931
$ws.Range("L51").Value = 0.01329
$ws.Range("F52").Value = 0.00404
$ws.Range("G52").Value = 0.0118370784417931
$ws.Range("L52").Value = 0.01329
$ws.Range("G53").Value = 0.517702347635114
$ws.Range("G54").Value = 0.517702347635114
$ws.Range("G56").Value = 0.532748510295894
$ws.Range("G57").Value = 0.532748510295894
$ws.Range("G66").Value = 1152.14464308153
$ws.Range("H66").Value = 13516.3435359425
$ws.Range("G67").Value = 1152.14464308153
$ws.Range("H67").Value = 13516.3435359425
$ws.Range("G68").Value = 1152.14464308153
$ws.Range("H68").Value = 13516.3435359425
$ws.Range("G69").Value = 1152.14464308153
$ws.Range("H69").Value = 13516.3435359425
$ws.Range("F71").Value = 0.00574
$ws.Range("G71").Value = 0.0132962454271556
$ws.Range("L71").Value = 0.01365
$ws.Range("M71").Value = 0.02337
$ws.Range("F72").Value = 0.00574
$ws.Range("G72").Value = 0.0132962454271556
$ws.Range("L72").Value = 0.01365
$ws.Range("M72").Value = 0.02337
$ws.Range("G73").Value = 0.43898706736516
$ws.Range("G74").Value = 0.43898706736516
$ws.Range("G76").Value = 0.45390986105187
$ws.Range("L76").Value = 0.14125
$ws.Range("G77").Value = 0.45390986105187
$ws.Range("L77").Value = 0.14125
$ws.Range("G86").Value = 1115.46043255521
$ws.Range("H86").Value = 13516.3435359425
$ws.Range("G87").Value = 1115.46043255521
$ws.Range("H87").Value = 13516.3435359425
$ws.Range("G88").Value = 1115.46043255521
$ws.Range("H88").Value = 13516.3435359425
$ws.Range("G89").Value = 1115.46043255521
$ws.Range("H89").Value = 13516.3435359425
$ws.Range("F91").Value = 0.00644
$ws.Range("G91").Value = 0.0110229414722533
$ws.Range("F92").Value = 0.00644
$ws.Range("G92").Value = 0.0110229414722533
$ws.Range("G93").Value = 0.497847775580576
$ws.Range("G94").Value = 0.497847775580576
$ws.Range("G96").Value = 0.5127326580096579
$ws.Range("L96").Value = 0.16625
$ws.Range("M96").Value = 1.06338
$ws.Range("G97").Value = 0.5127326580096579
$ws.Range("L97").Value = 0.16625
$ws.Range("M97").Value = 1.06338

# --- Add new rows 102-121 for year range 2019 - 2023 ---
# Row 102
$ws.Range("A102").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B102").Value = "ASPM"
$ws.Range("C102").Value = "B"
$ws.Range("D102").Value = "2019 - 2023"
$ws.Range("E102").Value = "RepSite"
$ws.Range("F102").Value = 0.441
$ws.Range("G102").Value = 0.3978
$ws.Range("H102").Value = 0.486
$ws.Range("I102").Value = 0.486
$ws.Range("L102").Value = 0.441
$ws.Range("M102").Value = 0.4741
$ws.Range("N102").Value = 0.486
$ws.Range("O102").Value = 1821031.91
$ws.Range("P102").Value = 5545591.67
$ws.Range("Q102").Value = "Manawatu District"
$ws.Range("R102").Value = "Manawatū"
$ws.Range("S102").Value = "Oroua"
$ws.Range("T102").Value = "Mana_12d"
# Row 103
$ws.Range("A103").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B103").Value = "Visual Clarity (Sediment class 3)"
$ws.Range("C103").Value = "D"
$ws.Range("D103").Value = "2019 - 2023"
$ws.Range("E103").Value = "RepSite"
$ws.Range("F103").Value = 1
$ws.Range("G103").Value = 1.58533333333333
$ws.Range("H103").Value = 5.1
$ws.Range("I103").Value = 4.95
$ws.Range("L103").Value = 1.05
$ws.Range("M103").Value = 3.1
$ws.Range("N103").Value = 3.87
$ws.Range("O103").Value = 1821031.91
$ws.Range("P103").Value = 5545591.67
$ws.Range("Q103").Value = "Manawatu District"
$ws.Range("R103").Value = "Manawatū"
$ws.Range("S103").Value = "Oroua"
$ws.Range("T103").Value = "Mana_12d"
$ws.Range("U103").Value = "m"
# Row 104
$ws.Range("A104").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B104").Value = "DRP (95th Percentile)"
$ws.Range("C104").Value = "C"
$ws.Range("D104").Value = "2019 - 2023"
$ws.Range("E104").Value = "RepSite"
$ws.Range("F104").Value = 0.012
$ws.Range("G104").Value = 0.0142881355932203
$ws.Range("H104").Value = 0.077
$ws.Range("I104").Value = 0.0343
$ws.Range("L104").Value = 0.008
$ws.Range("M104").Value = 0.019
$ws.Range("N104").Value = 0.02578
$ws.Range("O104").Value = 1821031.91
$ws.Range("P104").Value = 5545591.67
$ws.Range("Q104").Value = "Manawatu District"
$ws.Range("R104").Value = "Manawatū"
$ws.Range("S104").Value = "Oroua"
$ws.Range("T104").Value = "Mana_12d"
$ws.Range("U104").Value = "mg/L"
# Row 105
$ws.Range("A105").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B105").Value = "DRP (Median)"
$ws.Range("C105").Value = "C"
$ws.Range("D105").Value = "2019 - 2023"
$ws.Range("E105").Value = "RepSite"
$ws.Range("F105").Value = 0.012
$ws.Range("G105").Value = 0.0142881355932203
$ws.Range("H105").Value = 0.077
$ws.Range("I105").Value = 0.0343
$ws.Range("L105").Value = 0.008
$ws.Range("M105").Value = 0.019
$ws.Range("N105").Value = 0.02578
$ws.Range("O105").Value = 1821031.91
$ws.Range("P105").Value = 5545591.67
$ws.Range("Q105").Value = "Manawatu District"
$ws.Range("R105").Value = "Manawatū"
$ws.Range("S105").Value = "Oroua"
$ws.Range("T105").Value = "Mana_12d"
$ws.Range("U105").Value = "mg/L"
# Row 106
$ws.Range("A106").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B106").Value = "E coli (>260)"
$ws.Range("C106").Value = "C"
$ws.Range("D106").Value = "2019 - 2023"
$ws.Range("E106").Value = "RepSite"
$ws.Range("F106").Value = 146
$ws.Range("G106").Value = 965.745597511157
$ws.Range("H106").Value = 13516.3435359425
$ws.Range("I106").Value = 5807.2
$ws.Range("J106").Value = 22.4137931034483
$ws.Range("K106").Value = 31.0344827586207
$ws.Range("L106").Value = 100
$ws.Range("M106").Value = 964
$ws.Range("N106").Value = 3870.16
$ws.Range("O106").Value = 1821031.91
$ws.Range("P106").Value = 5545591.67
$ws.Range("Q106").Value = "Manawatu District"
$ws.Range("R106").Value = "Manawatū"
$ws.Range("S106").Value = "Oroua"
$ws.Range("T106").Value = "Mana_12d"
$ws.Range("U106").Value = "% exceedances over 260/100 mL"
# Row 107
$ws.Range("A107").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B107").Value = "E coli (>540)"
$ws.Range("C107").Value = "D"
$ws.Range("D107").Value = "2019 - 2023"
$ws.Range("E107").Value = "RepSite"
$ws.Range("F107").Value = 146
$ws.Range("G107").Value = 965.745597511157
$ws.Range("H107").Value = 13516.3435359425
$ws.Range("I107").Value = 5807.2
$ws.Range("J107").Value = 22.4137931034483
$ws.Range("K107").Value = 31.0344827586207
$ws.Range("L107").Value = 100
$ws.Range("M107").Value = 964
$ws.Range("N107").Value = 3870.16
$ws.Range("O107").Value = 1821031.91
$ws.Range("P107").Value = 5545591.67
$ws.Range("Q107").Value = "Manawatu District"
$ws.Range("R107").Value = "Manawatū"
$ws.Range("S107").Value = "Oroua"
$ws.Range("T107").Value = "Mana_12d"
$ws.Range("U107").Value = "% exceedances over 540/100 mL"
# Row 108
$ws.Range("A108").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B108").Value = "E coli (Median)"
$ws.Range("C108").Value = "D"
$ws.Range("D108").Value = "2019 - 2023"
$ws.Range("E108").Value = "RepSite"
$ws.Range("F108").Value = 146
$ws.Range("G108").Value = 965.745597511157
$ws.Range("H108").Value = 13516.3435359425
$ws.Range("I108").Value = 5807.2
$ws.Range("J108").Value = 22.4137931034483
$ws.Range("K108").Value = 31.0344827586207
$ws.Range("L108").Value = 100
$ws.Range("M108").Value = 964
$ws.Range("N108").Value = 3870.16
$ws.Range("O108").Value = 1821031.91
$ws.Range("P108").Value = 5545591.67
$ws.Range("Q108").Value = "Manawatu District"
$ws.Range("R108").Value = "Manawatū"
$ws.Range("S108").Value = "Oroua"
$ws.Range("T108").Value = "Mana_12d"
$ws.Range("U108").Value = "E. coli/100 mL"
# Row 109
$ws.Range("A109").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B109").Value = "E coli (95th Percentile)"
$ws.Range("C109").Value = "E"
$ws.Range("D109").Value = "2019 - 2023"
$ws.Range("E109").Value = "RepSite"
$ws.Range("F109").Value = 146
$ws.Range("G109").Value = 965.745597511157
$ws.Range("H109").Value = 13516.3435359425
$ws.Range("I109").Value = 5807.2
$ws.Range("J109").Value = 22.4137931034483
$ws.Range("K109").Value = 31.0344827586207
$ws.Range("L109").Value = 100
$ws.Range("M109").Value = 964
$ws.Range("N109").Value = 3870.16
$ws.Range("O109").Value = 1821031.91
$ws.Range("P109").Value = 5545591.67
$ws.Range("Q109").Value = "Manawatu District"
$ws.Range("R109").Value = "Manawatū"
$ws.Range("S109").Value = "Oroua"
$ws.Range("T109").Value = "Mana_12d"
$ws.Range("U109").Value = "E. coli/100 mL"
# Row 110
$ws.Range("A110").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B110").Value = "MCI"
$ws.Range("C110").Value = "D"
$ws.Range("D110").Value = "2019 - 2023"
$ws.Range("E110").Value = "RepSite"
$ws.Range("F110").Value = 89
$ws.Range("G110").Value = 93.774
$ws.Range("H110").Value = 111.11
$ws.Range("I110").Value = 111.11
$ws.Range("L110").Value = 89
$ws.Range("M110").Value = 106.553
$ws.Range("N110").Value = 111.11
$ws.Range("O110").Value = 1821031.91
$ws.Range("P110").Value = 5545591.67
$ws.Range("Q110").Value = "Manawatu District"
$ws.Range("R110").Value = "Manawatū"
$ws.Range("S110").Value = "Oroua"
$ws.Range("T110").Value = "Mana_12d"
# Row 111
$ws.Range("A111").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B111").Value = "Ammoniacal-N (95th Percentile)"
$ws.Range("C111").Value = "A"
$ws.Range("D111").Value = "2019 - 2023"
$ws.Range("E111").Value = "RepSite"
$ws.Range("F111").Value = 0.00717
$ws.Range("G111").Value = 0.0104006432611043
$ws.Range("H111").Value = 0.0499738233063448
$ws.Range("I111").Value = 0.03486
$ws.Range("L111").Value = 0.01159
$ws.Range("M111").Value = 0.01826
$ws.Range("N111").Value = 0.02758
$ws.Range("O111").Value = 1821031.91
$ws.Range("P111").Value = 5545591.67
$ws.Range("Q111").Value = "Manawatu District"
$ws.Range("R111").Value = "Manawatū"
$ws.Range("S111").Value = "Oroua"
$ws.Range("T111").Value = "Mana_12d"
$ws.Range("U111").Value = "mg NH4-N/L"
# Row 112
$ws.Range("A112").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B112").Value = "Ammoniacal-N (Median)"
$ws.Range("C112").Value = "A"
$ws.Range("D112").Value = "2019 - 2023"
$ws.Range("E112").Value = "RepSite"
$ws.Range("F112").Value = 0.00717
$ws.Range("G112").Value = 0.0104006432611043
$ws.Range("H112").Value = 0.0499738233063448
$ws.Range("I112").Value = 0.03486
$ws.Range("L112").Value = 0.01159
$ws.Range("M112").Value = 0.01826
$ws.Range("N112").Value = 0.02758
$ws.Range("O112").Value = 1821031.91
$ws.Range("P112").Value = 5545591.67
$ws.Range("Q112").Value = "Manawatu District"
$ws.Range("R112").Value = "Manawatū"
$ws.Range("S112").Value = "Oroua"
$ws.Range("T112").Value = "Mana_12d"
$ws.Range("U112").Value = "mg NH4-N/L"
# Row 113
$ws.Range("A113").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B113").Value = "Nitrate-N (95th Percentile)"
$ws.Range("C113").Value = "A"
$ws.Range("D113").Value = "2019 - 2023"
$ws.Range("E113").Value = "RepSite"
$ws.Range("F113").Value = 0.399
$ws.Range("G113").Value = 0.519494423452092
$ws.Range("H113").Value = 1.65
$ws.Range("I113").Value = 1.4555
$ws.Range("L113").Value = 0.181
$ws.Range("M113").Value = 1.06
$ws.Range("N113").Value = 1.2778
$ws.Range("O113").Value = 1821031.91
$ws.Range("P113").Value = 5545591.67
$ws.Range("Q113").Value = "Manawatu District"
$ws.Range("R113").Value = "Manawatū"
$ws.Range("S113").Value = "Oroua"
$ws.Range("T113").Value = "Mana_12d"
$ws.Range("U113").Value = "mg NO3-N/L"
# Row 114
$ws.Range("A114").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B114").Value = "Nitrate-N (Median)"
$ws.Range("C114").Value = "A"
$ws.Range("D114").Value = "2019 - 2023"
$ws.Range("E114").Value = "RepSite"
$ws.Range("F114").Value = 0.399
$ws.Range("G114").Value = 0.519494423452092
$ws.Range("H114").Value = 1.65
$ws.Range("I114").Value = 1.4555
$ws.Range("L114").Value = 0.181
$ws.Range("M114").Value = 1.06
$ws.Range("N114").Value = 1.2778
$ws.Range("O114").Value = 1821031.91
$ws.Range("P114").Value = 5545591.67
$ws.Range("Q114").Value = "Manawatu District"
$ws.Range("R114").Value = "Manawatū"
$ws.Range("S114").Value = "Oroua"
$ws.Range("T114").Value = "Mana_12d"
$ws.Range("U114").Value = "mg NO3-N/L"
# Row 115
$ws.Range("A115").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B115").Value = "QMCI"
$ws.Range("C115").Value = "B"
$ws.Range("D115").Value = "2019 - 2023"
$ws.Range("E115").Value = "RepSite"
$ws.Range("F115").Value = 5.787
$ws.Range("G115").Value = 5.5764
$ws.Range("H115").Value = 6.81
$ws.Range("I115").Value = 6.81
$ws.Range("L115").Value = 5.787
$ws.Range("M115").Value = 6.50025
$ws.Range("N115").Value = 6.81
$ws.Range("O115").Value = 1821031.91
$ws.Range("P115").Value = 5545591.67
$ws.Range("Q115").Value = "Manawatu District"
$ws.Range("R115").Value = "Manawatū"
$ws.Range("S115").Value = "Oroua"
$ws.Range("T115").Value = "Mana_12d"
# Row 116
$ws.Range("A116").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B116").Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Range("D116").Value = "2019 - 2023"
$ws.Range("E116").Value = "RepSite"
$ws.Range("F116").Value = 0.417
$ws.Range("G116").Value = 0.533127019738308
$ws.Range("H116").Value = 1.675
$ws.Range("I116").Value = 1.46665
$ws.Range("L116").Value = 0.20025
$ws.Range("M116").Value = 1.07578
$ws.Range("N116").Value = 1.30216
$ws.Range("O116").Value = 1821031.91
$ws.Range("P116").Value = 5545591.67
$ws.Range("Q116").Value = "Manawatu District"
$ws.Range("R116").Value = "Manawatū"
$ws.Range("S116").Value = "Oroua"
$ws.Range("T116").Value = "Mana_12d"
$ws.Range("U116").Value = "g/m3"
# Row 117
$ws.Range("A117").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B117").Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Range("D117").Value = "2019 - 2023"
$ws.Range("E117").Value = "RepSite"
$ws.Range("F117").Value = 0.417
$ws.Range("G117").Value = 0.533127019738308
$ws.Range("H117").Value = 1.675
$ws.Range("I117").Value = 1.46665
$ws.Range("L117").Value = 0.20025
$ws.Range("M117").Value = 1.07578
$ws.Range("N117").Value = 1.30216
$ws.Range("O117").Value = 1821031.91
$ws.Range("P117").Value = 5545591.67
$ws.Range("Q117").Value = "Manawatu District"
$ws.Range("R117").Value = "Manawatū"
$ws.Range("S117").Value = "Oroua"
$ws.Range("T117").Value = "Mana_12d"
$ws.Range("U117").Value = "g/m3"
# Row 118
$ws.Range("A118").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B118").Value = "Total Nitrogen (95th Percentile)"
$ws.Range("D118").Value = "2019 - 2023"
$ws.Range("E118").Value = "RepSite"
$ws.Range("F118").Value = 0.79
$ws.Range("G118").Value = 0.902203389830509
$ws.Range("H118").Value = 2.88
$ws.Range("I118").Value = 1.991
$ws.Range("L118").Value = 0.5649999999999999
$ws.Range("M118").Value = 1.3888
$ws.Range("N118").Value = 1.909
$ws.Range("O118").Value = 1821031.91
$ws.Range("P118").Value = 5545591.67
$ws.Range("Q118").Value = "Manawatu District"
$ws.Range("R118").Value = "Manawatū"
$ws.Range("S118").Value = "Oroua"
$ws.Range("T118").Value = "Mana_12d"
$ws.Range("U118").Value = "g/m3"
# Row 119
$ws.Range("A119").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B119").Value = "Total Nitrogen (Median)"
$ws.Range("D119").Value = "2019 - 2023"
$ws.Range("E119").Value = "RepSite"
$ws.Range("F119").Value = 0.79
$ws.Range("G119").Value = 0.902203389830509
$ws.Range("H119").Value = 2.88
$ws.Range("I119").Value = 1.991
$ws.Range("L119").Value = 0.5649999999999999
$ws.Range("M119").Value = 1.3888
$ws.Range("N119").Value = 1.909
$ws.Range("O119").Value = 1821031.91
$ws.Range("P119").Value = 5545591.67
$ws.Range("Q119").Value = "Manawatu District"
$ws.Range("R119").Value = "Manawatū"
$ws.Range("S119").Value = "Oroua"
$ws.Range("T119").Value = "Mana_12d"
$ws.Range("U119").Value = "g/m3"
# Row 120
$ws.Range("A120").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B120").Value = "Total Phosphorus (95th Percentile)"
$ws.Range("D120").Value = "2019 - 2023"
$ws.Range("E120").Value = "RepSite"
$ws.Range("F120").Value = 0.025
$ws.Range("G120").Value = 0.0654237288135593
$ws.Range("H120").Value = 0.5590000000000001
$ws.Range("I120").Value = 0.25535
$ws.Range("L120").Value = 0.0205
$ws.Range("M120").Value = 0.10576
$ws.Range("N120").Value = 0.17534
$ws.Range("O120").Value = 1821031.91
$ws.Range("P120").Value = 5545591.67
$ws.Range("Q120").Value = "Manawatu District"
$ws.Range("R120").Value = "Manawatū"
$ws.Range("S120").Value = "Oroua"
$ws.Range("T120").Value = "Mana_12d"
$ws.Range("U120").Value = "g/m3"
# Row 121
$ws.Range("A121").Value = "Kiwitea at Kimbolton Rd"
$ws.Range("B121").Value = "Total Phosphorus (Median)"
$ws.Range("D121").Value = "2019 - 2023"
$ws.Range("E121").Value = "RepSite"
$ws.Range("F121").Value = 0.025
$ws.Range("G121").Value = 0.0654237288135593
$ws.Range("H121").Value = 0.5590000000000001
$ws.Range("I121").Value = 0.25535
$ws.Range("L121").Value = 0.0205
$ws.Range("M121").Value = 0.10576
$ws.Range("N121").Value = 0.17534
$ws.Range("O121").Value = 1821031.91
$ws.Range("P121").Value = 5545591.67
$ws.Range("Q121").Value = "Manawatu District"
$ws.Range("R121").Value = "Manawatū"
$ws.Range("S121").Value = "Oroua"
$ws.Range("T121").Value = "Mana_12d"
$ws.Range("U121").Value = "g/m3"

# --- Set explicit empty-string cells (trailing apostrophe -> empty text, not a blank cell) ---
$ws.Range("J102").Value = "'"
$ws.Range("J102").Style = "Normal"
$ws.Range("K102").Value = "'"
$ws.Range("K102").Style = "Normal"
$ws.Range("U102").Value = "'"
$ws.Range("U102").Style = "Normal"
$ws.Range("J103").Value = "'"
$ws.Range("J103").Style = "Normal"
$ws.Range("K103").Value = "'"
$ws.Range("K103").Style = "Normal"
$ws.Range("J104").Value = "'"
$ws.Range("J104").Style = "Normal"
$ws.Range("K104").Value = "'"
$ws.Range("K104").Style = "Normal"
$ws.Range("J105").Value = "'"
$ws.Range("J105").Style = "Normal"
$ws.Range("K105").Value = "'"
$ws.Range("K105").Style = "Normal"
$ws.Range("J110").Value = "'"
$ws.Range("J110").Style = "Normal"
$ws.Range("K110").Value = "'"
$ws.Range("K110").Style = "Normal"
$ws.Range("U110").Value = "'"
$ws.Range("U110").Style = "Normal"
$ws.Range("J111").Value = "'"
$ws.Range("J111").Style = "Normal"
$ws.Range("K111").Value = "'"
$ws.Range("K111").Style = "Normal"
$ws.Range("J112").Value = "'"
$ws.Range("J112").Style = "Normal"
$ws.Range("K112").Value = "'"
$ws.Range("K112").Style = "Normal"
$ws.Range("J113").Value = "'"
$ws.Range("J113").Style = "Normal"
$ws.Range("K113").Value = "'"
$ws.Range("K113").Style = "Normal"
$ws.Range("J114").Value = "'"
$ws.Range("J114").Style = "Normal"
$ws.Range("K114").Value = "'"
$ws.Range("K114").Style = "Normal"
$ws.Range("J115").Value = "'"
$ws.Range("J115").Style = "Normal"
$ws.Range("K115").Value = "'"
$ws.Range("K115").Style = "Normal"
$ws.Range("U115").Value = "'"
$ws.Range("U115").Style = "Normal"
$ws.Range("C116").Value = "'"
$ws.Range("C116").Style = "Normal"
$ws.Range("J116").Value = "'"
$ws.Range("J116").Style = "Normal"
$ws.Range("K116").Value = "'"
$ws.Range("K116").Style = "Normal"
$ws.Range("C117").Value = "'"
$ws.Range("C117").Style = "Normal"
$ws.Range("J117").Value = "'"
$ws.Range("J117").Style = "Normal"
$ws.Range("K117").Value = "'"
$ws.Range("K117").Style = "Normal"
$ws.Range("C118").Value = "'"
$ws.Range("C118").Style = "Normal"
$ws.Range("J118").Value = "'"
$ws.Range("J118").Style = "Normal"
$ws.Range("K118").Value = "'"
$ws.Range("K118").Style = "Normal"
$ws.Range("C119").Value = "'"
$ws.Range("C119").Style = "Normal"
$ws.Range("J119").Value = "'"
$ws.Range("J119").Style = "Normal"
$ws.Range("K119").Value = "'"
$ws.Range("K119").Style = "Normal"
$ws.Range("C120").Value = "'"
$ws.Range("C120").Style = "Normal"
$ws.Range("J120").Value = "'"
$ws.Range("J120").Style = "Normal"
$ws.Range("K120").Value = "'"
$ws.Range("K120").Style = "Normal"
$ws.Range("C121").Value = "'"
$ws.Range("C121").Style = "Normal"
$ws.Range("J121").Value = "'"
$ws.Range("J121").Style = "Normal"
$ws.Range("K121").Value = "'"
$ws.Range("K121").Style = "Normal"
